# Updated cryptos list on Fri Jun  9 20:21:00 UTC 2023 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values
# for the cryptocurrency rows in the worksheet. D values that look like
# plain numbers are prefixed with a leading apostrophe so Excel keeps
# them stored as text (matching the source data's text format), the
# same way the original report text values (e.g. thousand-dot prices)
# are preserved as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.466.87'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.834.29'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'260.24"
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("D7").Value = "'0.5381"
$ws.Range("E7").Value = '  +2.61%  '
$ws.Range("D8").Value = "'0.3019"
$ws.Range("E8").Value = '  -6.83%  '
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").Value = "'17.49"
$ws.Range("E10").Value = '  -7.35%  '
$ws.Range("D11").Value = "'0.7358"
$ws.Range("E11").Value = '  -5.75%  '
$ws.Range("D12").Value = '1.845.04'
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").Value = "'0.07219"
$ws.Range("E13").Value = '  -7.27%  '
$ws.Range("D14").Value = "'88.92"
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = "'4.953"
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = "'13.80"
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = "'0.000007868"
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("D20").Value = '26.491.57'
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").Value = '2.080.83'
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").Value = "'4.559"
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").Value = "'9.219"
$ws.Range("E24").Value = '  -2.82%  '
$ws.Range("D25").Value = "'142.50"
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = "'2.198"
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").Value = "'1.684"
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = "'16.95"
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("D29").Value = "'110.32"
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").Value = "'4.204"
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = "'0.08790"
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("D32").Value = "'3.989"
$ws.Range("E32").Value = '  -3.21%  '
$ws.Range("D33").Value = "'0.04786"
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").Value = "'2.938"
$ws.Range("E34").Value = '  +2.05%  '
$ws.Range("D35").Value = "'0.7265"
$ws.Range("E35").Value = '  +0.60%  '
$ws.Range("D36").Value = "'1.127"
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").Value = "'3.087"
$ws.Range("E37").Value = '  -0.64%  '
$ws.Range("D38").Value = "'2.270"
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("D39").Value = "'0.01701"
$ws.Range("E39").Value = '  -5.09%  '
$ws.Range("D40").Value = "'0.4705"
$ws.Range("E40").Value = '  -3.32%  '
$ws.Range("D41").Value = "'0.9047"
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = "'107.63"
$ws.Range("E42").Value = '  -2.87%  '
$ws.Range("D43").Value = "'5.874"
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").Value = "'7.325"
$ws.Range("D46").Value = "'8.933"
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("D47").Value = "'0.1228"
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").Value = "'0.4043"
$ws.Range("E48").Value = '  -4.05%  '
$ws.Range("D49").Value = "'0.05787"
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("D50").Value = "'34.62"
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").Value = "'0.8871"
$ws.Range("E51").Value = '  -0.26%  '
